# CouponLoyaltyRTM.xlsx update: refresh RTM entries (sentence #s, priorities,
# and two use-case descriptions) per the revised spec.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM")

# Reworded / renamed use cases (kept in this order so the shared-strings
# table lines up with the authored edit).
$ws.Range("C4").Value = "A loyalty program system that uses the user's phone number, or email to create a unique ID. "
$ws.Range("C10").Value = "A way for users to generate reports and filter on date, product, etc."
$ws.Range("D10").Value = "UC8_Generate_Report"
$ws.Range("D4").Value = "UC2_Loyalty_Tracking_System:"

# Row 3 - UC1_User_Login
$ws.Range("B3").Value = 3
$ws.Range("F3").Value = 2

# Row 4 - UC2_Loyalty_Tracking_System
$ws.Range("F4").Value = 1

# Row 5 - UC3_View_Sales_Report
$ws.Range("F5").Value = 3

# Row 6 - UC4_Customer_Absent_Report
$ws.Range("F6").Value = 6

# Row 7 - UC5_Coupon_Ratio_Report
$ws.Range("F7").Value = 8

# Row 8 - UC6_Create_Coupon
$ws.Range("B8").Value = 4
$ws.Range("F8").Value = 9

# Row 9 - UC7_View_Available_Coupons
$ws.Range("B9").Value = 4
$ws.Range("F9").Value = 7

# Row 10 - UC8_Generate_Report priority
$ws.Range("F10").Value = 4

# Row 11 - UC9_Frequent_Buyer_Report
$ws.Range("F11").Value = 5

# Row 12 - UC10_Change_Password
$ws.Range("F12").Value = 10

# Match the author's final selection on the RTM sheet
$ws.Range("D7").Select()
